$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# 2. Convert the data range into an Excel Table (ListObject)
$range = $ws.Range("A1:U86")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# 3. Freeze the header row (pane split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
